$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("journalVoucherDetails")

# Remove the now-unused "voucherDate" column (old column B). Every column to
# its right (voucherType, accountCode1, accountCode2, department, function,
# accountCode3) shifts one place to the left, taking its data/format with it.
$ws.Columns.Item(2).Delete()

# Rename two of the existing test rows (the underlying accountCode columns
# already carry the right values/format after the column shift above).
$ws.Range("A4").Value = "voucherBillPayment"
$ws.Range("A5").Value = "budgetCheckWithSubledger"

# New row 7: budgetCheckWithOutSubledger
$ws.Range("A7").Value = "budgetCheckWithOutSubledger"
$ws.Range("B7").Value = "Expense"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "1100101"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3501003"
$ws.Range("E7").Value = "ENGINEERING"
$ws.Range("F7").Value = "Water Supply"

# New row 8: voucherWithOutSubledger
$ws.Range("A8").Value = "voucherWithOutSubledger"
$ws.Range("B8").Value = "Expense"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "1100101"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3501003"
$ws.Range("E8").Value = "ENGINEERING"
$ws.Range("F8").Value = "Water Supply"

# New row 9: voucherWithSubledger
$ws.Range("A9").Value = "voucherWithSubledger"
$ws.Range("B9").Value = "General"
$ws.Range("C9").NumberFormat = "General"
$ws.Range("C9").Value = 2101001
$ws.Range("C9").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Value = 3501003
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").Value = "PUBLIC HEALTH AND SANITATION"
$ws.Range("F9").Value = "Public Health"

# Match the author's final selection/active cell on this sheet.
$ws.Range("A4").Select()
